# Autocad/AddPoindWithCalcZ/check.xlsx
#
# Commit: Added "TEXT" to function BRfromMtext: TEXT needs to be in layers:
# 1506, 1507. br will be created in layers M1506, M1507 accordingly.
#
# The accompanying check.xlsx recomputes the sample point set used to
# validate the new-Z interpolation math. Update the "start"/"end"/"new"
# sample coordinates on the "Sheet1 (2)" worksheet; the dependent formulas
# (F2/G2, row 7 and row 8) recalculate automatically. Also restore the
# last-used selection on that sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sheet1 (2)")

# row 2 ("start")
$ws.Range("B2").Value = 204091.13
$ws.Range("C2").Value = 741349.81
$ws.Range("D2").Value = 181.33

# row 3 ("end")
$ws.Range("B3").Value = 204095.88
$ws.Range("C3").Value = 741350.14
$ws.Range("D3").Value = 180.71

# row 4 ("new")
$ws.Range("B4").Value = 204093.24
$ws.Range("C4").Value = 741350.36

# Leave the last active selection on C4, matching the saved view state.
$ws.Activate() | Out-Null
$ws.Range("C4").Select() | Out-Null
